# "Generate Report for Handoff"
#
# The handoff/generate-report run for the 124a4909-3162-4fd1-98d4-c7460ba0638b
# file produced new timestamps. Update the three cells that record them:
#   - Overview!G6  "Latest HO Xliff Generate Date"
#   - zh-cn!H6     "Latest Handoff Datetime"
#   - de-de!H6     "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(6, 7).Value = "2016-08-31 10:46:38"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(6, 8).Value = "2016-08-31 10:46:34"

$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(6, 8).Value = "2016-08-31 10:46:38"
